# Exercice 4 - Instructions - Position absolute et fixed
# Applies the content edit described in the commit diff:
#  1. Center the title paragraph.
#  2. Rewrite the "Objectif" explanation paragraph (new wording + a manual
#     line break) and turn the old "exigences" sentence into four bulleted
#     requirements.
#  3. Give the illustration a picture border/fill (closest achievable via
#     the InlineShape.Line / .Fill object model) and bump its editId.
#  4. Drop the trailing empty paragraph at the end of the body.

$d = $word.ActiveDocument

# 1) Center-align the big heading ("Exercice #4 - ...").
$d.Paragraphs(1).Alignment = 1   # wdAlignParagraphCenter -> <w:jc w:val="center"/>

# 2) Replace the "Utiliser les positionnements ..." paragraph's text.
$targetPara = $d.Paragraphs(4)
$rng = $targetPara.Range
$rng.MoveEnd(1, -1) | Out-Null            # keep the paragraph mark out of the range
$vt = [char]11                             # manual line break (-> <w:br/>)
$rng.Text = "Ce numéro est composé de 4 petits énoncés" + $vt + "pour lesquelles vous devez positionner des rectangles tel que décrit dans le HTML. "

# 3) Turn the old "exigences" sentence into four bulleted list items right
#    after that paragraph, all sharing the same bullet list (same numId).
$targetPara.Range.InsertParagraphAfter()
$item1 = $d.Paragraphs(5)
$item1.Range.Text = "Vous ne devez pas modifier le HTML"
$item1.Range.ListFormat.ApplyBulletDefault()
$bulletTemplate = $item1.Range.ListFormat.ListTemplate

$item1.Range.InsertParagraphAfter()
$item2 = $d.Paragraphs(6)
$item2.Range.Text = "Vous pouvez uniquement modifier les règles CSS présentes dans l'élément style du head"
$item2.Range.ListFormat.ApplyListTemplateWithLevel($bulletTemplate, $true)

$item2.Range.InsertParagraphAfter()
$item3 = $d.Paragraphs(7)
$item3.Range.Text = "Vous devez uniquement utilisez les positionnements absolute et fixed"
$item3.Range.ListFormat.ApplyListTemplateWithLevel($bulletTemplate, $true)

$item3.Range.InsertParagraphAfter()
$item4 = $d.Paragraphs(8)
$item4.Range.Text = "Vous pouvez utiliser les propriétés : left, right, top, bottom"
$item4.Range.ListFormat.ApplyListTemplateWithLevel($bulletTemplate, $true)

# 4) Apply a picture border + subtle dark fill to the illustration, mirroring
#    the "Picture Border" formatting captured in the diff.
$shp = $d.InlineShapes(1)
$shp.Fill.Visible = $true
$shp.Fill.ForeColor.RGB = 0
$shp.Line.Visible = $true
$shp.Line.Weight = 0.25
$shp.Line.ForeColor.RGB = 0
$shp.Line.Style = 1
$shp.Line.DashStyle = 1

# 5) Remove the now-superfluous empty paragraph at the very end of the body.
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$secondLastPara = $d.Paragraphs($d.Paragraphs.Count - 1)
if ($lastPara.Range.Text -eq [char]13 -and $secondLastPara.Range.InlineShapes.Count -eq 0) {
    $lastPara.Range.Delete()
}
